$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.937.79'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').Value = '1.833.82'
$ws.Range('E3').Value = '  -1.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.52'
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6903'
$ws.Range('E6').Value = '  -2.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9997'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07705'
$ws.Range('E8').Value = '  -2.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3053'
$ws.Range('E9').Value = '  -2.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.55'
$ws.Range('E10').Value = '  -4.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07802'
$ws.Range('E11').Value = '  -0.97%  '
$ws.Range('D12').Value = '1.839.73'
$ws.Range('E12').Value = '  -1.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.082'
$ws.Range('E13').Value = '  -2.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '90.50'
$ws.Range('E14').Value = '  -3.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6805'
$ws.Range('E15').Value = '  -3.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.439'
$ws.Range('E16').Value = '  -1.45%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008346'
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('D18').Value = '28.942.65'
$ws.Range('E18').Value = '  -1.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.68'
$ws.Range('E19').Value = '  -4.26%  '
$ws.Range('D20').Value = '2.083.54'
$ws.Range('E20').Value = '  -1.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.71'
$ws.Range('E21').Value = '  -3.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9993'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.476'
$ws.Range('E23').Value = '  -2.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9994'
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.29'
$ws.Range('E25').Value = '  +0.69%  '
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1471'
$ws.Range('E26').Value = '  -5.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.799'
$ws.Range('E27').Value = '  -2.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.22'
$ws.Range('E28').Value = '  -3.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.555'
$ws.Range('E29').Value = '  +3.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.217'
$ws.Range('E30').Value = '  -2.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.157'
$ws.Range('E31').Value = '  -2.47%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.179'
$ws.Range('E32').Value = '  -2.91%  '
$ws.Range('E33').Value = '  -3.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7649'
$ws.Range('E34').Value = '  +1.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.849'
$ws.Range('E35').Value = '  -2.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.147'
$ws.Range('E36').Value = '  -2.53%  '
$ws.Range('E37').Value = '  -1.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01847'
$ws.Range('E38').Value = '  -2.04%  '
$ws.Range('D39').Value = '1.234.94'
$ws.Range('E39').Value = '  -3.79%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.695'
$ws.Range('E40').Value = '  -2.45%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9217'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '108.54'
$ws.Range('E42').Value = '  -0.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.860'
$ws.Range('E43').Value = '  -2.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9992'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '9.570'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('E46').Value = '  -4.37%  '
$ws.Range('D47').Value = '1.982.69'
$ws.Range('E47').Value = '  -2.33%  '
$ws.Range('E48').Value = '  -0.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '64.23'
$ws.Range('E49').Value = '  -9.70%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.747'
$ws.Range('E50').Value = '  -2.90%  '
$ws.Range('E51').Value = '  -1.95%  '
